# Build site at 2023-04-12 14:53:07 UTC
# Updates the LOQ4058 course-syllabus sheet:
#  - "Objetivos:" (row 10) gets its real objectives text (it had been
#    holding the professor's name by mistake).
#  - A new row is inserted right under "Docentes responsaveis:" to hold
#    the professor's name on its own row.
#  - "Programa resumido:", "Programa:", "Metodo:", "Criterio:",
#    "Norma de recuperacao:" and "Bibliografia:" each receive their
#    real body text (previously they were holding placeholder/borrowed
#    values from neighbouring rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix the "Objetivos:" row (row 10): it currently holds the
#    professor's name; it should hold the real course-objectives text.
$ws.Range("B10").Value = "O principal objetivo é permitir ao aluno conhecer os princípios fundamentais das interfaces líquido/gás/sólido e aprender sobre os conceitos de adsorção em sólidos, a caracterização de superfícies porosas, tensão superficial e propriedades de sistemas coloidais e emulsões. Além disso, identificar e explorar as aplicações destes conceitos em diferentes processos da indústria química."
$ws.Range("C10").Value = "O principal objetivo é permitir ao aluno conhecer os princípios fundamentais das interfaces líquido/gás/sólido e aprender sobre os conceitos de adsorção em sólidos, a caracterização de superfícies porosas, tensão superficial e propriedades de sistemas coloidais e emulsões. Além disso, identificar e explorar as aplicações destes conceitos em diferentes processos da indústria química."

# 2) Insert a new row 13 (right after "Docentes responsaveis:", row 12)
#    to hold the professor's name that used to sit on row 10.
$ws.Rows.Item(13).Insert()

# The freshly-inserted row inherits row 12's formatting (label style in
# column A). Column A must stay empty/unstyled on the new row, so reset
# it back to the default "Normal" style and drop its content entirely.
$ws.Range("A13").Style = "Normal"
$ws.Range("A13").ClearContents()

# Columns B/C on the new row must use the normal body-text styles (the
# same ones already used one row below, on B14/C14) instead of row 12's
# label style, so copy just the formatting down before setting values.
$ws.Range("B14").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

$ws.Range("B13").Value = "1488970 - Marivone Nunho Sousa"
$ws.Range("C13").Value = "1488970 - Marivone Nunho Sousa"

# 3) "Programa resumido:" (now row 14) gets the real short-syllabus text.
$ws.Range("B14").Value = "Descrição de superfície e interface, termodinâmica das superfícies. Superfícies e forças. Sólidos iônicos e covalentes. Forças físicas e químicas de adsorção. Interface gás-solido e líquido-sólido. Caracterização de superfícies."
$ws.Range("C14").Value = "Descrição de superfície e interface, termodinâmica das superfícies. Superfícies e forças. Sólidos iônicos e covalentes. Forças físicas e químicas de adsorção. Interface gás-solido e líquido-sólido. Caracterização de superfícies."

# 4) "Programa:" (now row 16) gets the full syllabus text.
$ws.Range("B16").Value = "1)Isotermas de adsorção: Isotermas de Langmuir. Isotermas de Brunauer, Emmett e Teller. Métodos de determinação da área superficial de sólidos. Classificação quanto à porosidade. Métodos de determinação da porosidade de sólidos.2)Catálise de superfície. Interface sólido‐líquido. Mecanismos de catálise.3)Tensão superficial e interfacial. Equação de Laplace. Ângulo de contato. Ascensão e depressão capilar. Aplicações.4)Classificação das dispersões coloidais. Dupla camada elétrica: equação de Lippman e apresentação de modelos. Estabilidade e coagulação de dispersões coloidais. 5)Interações intermoleculares, dipolo-dipolo e de Van-der-Waals, ligação de hidrogênio e interações estabilizadoras em macromoléculas. 6)Estado coloidal. Colóides liofílicos e liofóbicos, hdrofílicos e hidrofóbicos. Obtenção de colóides. Propriedades cinéticas difusão, sedimentação, convecção. Propriedades óticas: espalhamento estático de luz, turbidez, espalhamento dinâmico da luz.Coagulação. Aplicações.7)Termodinâmica dos processos de transporte: difusão sedimentação e transporte através de membranas."
$ws.Range("C16").Value = "1)Isotermas de adsorção: Isotermas de Langmuir. Isotermas de Brunauer, Emmett e Teller. Métodos de determinação da área superficial de sólidos. Classificação quanto à porosidade. Métodos de determinação da porosidade de sólidos.2)Catálise de superfície. Interface sólido‐líquido. Mecanismos de catálise.3)Tensão superficial e interfacial. Equação de Laplace. Ângulo de contato. Ascensão e depressão capilar. Aplicações.4)Classificação das dispersões coloidais. Dupla camada elétrica: equação de Lippman e apresentação de modelos. Estabilidade e coagulação de dispersões coloidais. 5)Interações intermoleculares, dipolo-dipolo e de Van-der-Waals, ligação de hidrogênio e interações estabilizadoras em macromoléculas. 6)Estado coloidal. Colóides liofílicos e liofóbicos, hdrofílicos e hidrofóbicos. Obtenção de colóides. Propriedades cinéticas difusão, sedimentação, convecção. Propriedades óticas: espalhamento estático de luz, turbidez, espalhamento dinâmico da luz.Coagulação. Aplicações.7)Termodinâmica dos processos de transporte: difusão sedimentação e transporte através de membranas."

# 5) "Metodo:" (now row 19) gets the evaluation-method text.
$ws.Range("B19").Value = "Participação em sala de aula, preparação e apresentação de trabalhos e provas escritas."
$ws.Range("C19").Value = "Participação em sala de aula, preparação e apresentação de trabalhos e provas escritas."

# 6) "Criterio:" (now row 20) gets the final-grade formula text.
$ws.Range("B20").Value = "Média Final = (Prova1 + Prova2 + Nota de Trabalho) /3`nMédia final mínima de aprovação = 5,0"
$ws.Range("C20").Value = "Média Final = (Prova1 + Prova2 + Nota de Trabalho) /3`nMédia final mínima de aprovação = 5,0"

# 7) "Norma de recuperacao:" (now row 21) gets the make-up exam rule text.
$ws.Range("B21").Value = "(Prova escrita + Média Final)/2         Nota Final mínima para aprovação= 5,0"
$ws.Range("C21").Value = "(Prova escrita + Média Final)/2         Nota Final mínima para aprovação= 5,0"

# 8) "Bibliografia:" (now row 22) gets the real bibliography text.
$ws.Range("B22").Value = "1)MYERS, D. Surfaces, interfaces, and colloids: Principles and Applications, Second edition, Wiley-VCH, New York, 19982) BIRDI, K. S.; Surface and Colloid Chemistry, 1a ed., CRC Press LLC, New York, 1997.3) OSHIMA, H., Theory of colloid and interfacial electric phenomena. Interface Science and Technology Series, v. 12, Academic Press, Oxford, 2006.4) JACOB N. ISRAELCHVILI; Intermolecular and Surface Forces, 3r d Edition, New York, Academic, 2010.5) ADAMIAN, R. E ALMENDRA E.; Físico-Química – Uma Aplicação aos Materiais, 2002. 6) ADAMSON, A. Physical Chemistry of Surfaces (5th ed.). New York: John Wiley, 1990.7) SHAW, D. J. Introdução à Química dos Coloides e de Superfícies. São Paulo: Edgard Blücher, 1975. 185 pp.8)  REGALBUTO, J. Handbook of catalyst preparation. Taylor & Francis,2007"
$ws.Range("C22").Value = "1)MYERS, D. Surfaces, interfaces, and colloids: Principles and Applications, Second edition, Wiley-VCH, New York, 19982) BIRDI, K. S.; Surface and Colloid Chemistry, 1a ed., CRC Press LLC, New York, 1997.3) OSHIMA, H., Theory of colloid and interfacial electric phenomena. Interface Science and Technology Series, v. 12, Academic Press, Oxford, 2006.4) JACOB N. ISRAELCHVILI; Intermolecular and Surface Forces, 3r d Edition, New York, Academic, 2010.5) ADAMIAN, R. E ALMENDRA E.; Físico-Química – Uma Aplicação aos Materiais, 2002. 6) ADAMSON, A. Physical Chemistry of Surfaces (5th ed.). New York: John Wiley, 1990.7) SHAW, D. J. Introdução à Química dos Coloides e de Superfícies. São Paulo: Edgard Blücher, 1975. 185 pp.8)  REGALBUTO, J. Handbook of catalyst preparation. Taylor & Francis,2007"

Write-Host "edit applied"
